# "Add change log and roadmap": append a new entry "Dati integrativi AIRE"
# (code 21) to the "Tipi di mutazione" lookup table on Foglio1.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the formatting of the last existing data row (A21:B21) down onto the
# new row 22, then overwrite its values — this carries over the style
# (font/border/alignment) used by every other data row in the table.
$ws.Range("A21:B21").Copy()
$ws.Range("A22:B22").PasteSpecial(-4122)   # xlPasteFormats

$ws.Range("A22").Value = 21
$ws.Range("B22").Value = "Dati integrativi AIRE"

# Leave the selection where the author left it after adding the row.
$ws.Range("E20").Select()
